$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 48; this pushes the existing
# rows 48-161 down to 49-162 (matches the diff: every row from 49
# downward now carries what used to be one row above it, and a new
# row 162 appears carrying the former row 161's data).
$ws.Rows(48).Insert()

# Populate the newly inserted row 48 with the new record.
$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = 44519
$ws.Range("E48").Value = 10
$ws.Range("F48").Value = 100112017
$ws.Range("G48").Value = "Apio"
$ws.Range("H48").Value = "Americana (o)"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 45
$ws.Range("K48").Value = 12000
$ws.Range("L48").Value = 12000
$ws.Range("M48").Value = 12000
$ws.Range("N48").Value = "`$/docena de matas"
$ws.Range("O48").Value = "Región de Coquimbo"
$ws.Range("P48").Value = 2000
$ws.Range("Q48").Value = 6
$ws.Range("R48").Value = "Hortaliza"
